$wb = $excel.ActiveWorkbook

# Rename the first sheet ("Simply" -> "Headers") per commit diff
$ws = $wb.Worksheets.Item("Simply")
$ws.Name = "Headers"
$ws.Activate()

# B1 header changes from shared-string index 65 ("Manual_Test_Case_Folder")
# to 64 after four "...(Simply...)..." strings are deleted elsewhere in the
# shared-string table; net visible effect is the same text stays in B1.
$ws.Range("B1").Value = "Manual_Test_Case_Folder"

# D1:G1 header labels (Address1/City/State/Zip) are removed
$ws.Range("D1:G1").ClearContents()

# Rows 2-5: drop the per-row Simply-XX test case name (C) and the
# Address1/City/State/Zip/Date sample data (D:G)
$ws.Range("C2:G5").ClearContents()

# New selection/active cell recorded in the sheet view
$ws.Range("N26").Select()
